$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (Player, Position, Team) for rows 2-19
$data = @(
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Kyshawn George", "SG,SF", "Washington Wizards"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Amen Thompson", "SG,SF,PF", "Houston Rockets"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
